$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for several rows after repulling data
$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -17
$ws.Range("F10").Value = -9
$ws.Range("F14").Value = -8
